$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value changes -------------------------------------------------
# B3 / B4 get the same "PlantUML Complete" marker that B2 already has.
$ws.Range("B3").Value = "PlantUML Complete"
$ws.Range("B4").Value = "PlantUML Complete"

# C2 gets a brand new note.
$ws.Range("C2").Value = "Fix Activation Buttons"

# --- Fill colors ----------------------------------------------------------
# Green ("Green, Accent 6, Lighter 80%") for the "PlantUML Complete" column,
# gold ("Gold, Accent 4, Lighter 80%") for the new activation-buttons note.
$green = 14348258   # RGB(226,239,218)
$gold  = 13431551   # RGB(255,242,204)

$ws.Range("B2").Interior.Color = $green
$ws.Range("B3").Interior.Color = $green
$ws.Range("B4").Interior.Color = $green
$ws.Range("C2").Interior.Color = $gold

# --- Selection --------------------------------------------------------
$ws.Range("E6").Select()
